# Add a new menu item ("Kinder Creamy") to the Snacks category.
# This inserts a new row 49 (pushing the existing rows 49-61 down to 50-62),
# fills in the new row's values, and updates the hidden AutoFilter
# defined name range to account for the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49 (Snacks section), shifting rows 49:61 -> 50:62
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new item's data
$ws.Range("A49").Value = "Snacks"
$ws.Range("B49").Value = "Kinder Creamy"
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 25
$ws.Range("E49").Value = "Kinder Creamy.jpg"
$ws.Range("F49").Value = "Fast Food"

# The sheet now spans down to row 62; update the hidden _FilterDatabase
# defined name (tracks the last AutoFilter range) to match.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$F`$60"

# Scroll the view down so row 43 becomes the top visible row (matches the
# author's saved scroll position after adding the new row).
$excel.ActiveWindow.ScrollRow = 43
